# CPI2019.xlsx edit
# - Rename the "CPI Timeseries 2012 - 2019" sheet to "CPI2012-2019"
# - Fix up the Print_Titles defined name so it still points at the renamed sheet
#   (the _FilterDatabase name updates itself automatically on rename)
# - Make the renamed "CPI2012-2019" sheet the active sheet/tab, with its
#   frozen-pane selection moved to C101 (was A113)

$wb = $excel.ActiveWorkbook

$tsSheet = $wb.Worksheets.Item("CPI Timeseries 2012 - 2019")
$tsSheet.Name = "CPI2012-2019"

foreach ($n in $wb.Names) {
    if ($n.Name -eq "CPI2012-2019!Print_Titles") {
        $n.RefersTo = "='CPI2012-2019'!`$A:`$A,'CPI2012-2019'!`$3:`$3"
    }
}

$tsSheet.Activate()
$tsSheet.Range("C101").Select()
